# Update "correct numbers for D2-C2 scenarios with medium density" across
# the "Low traffic densit" and "Medium traffic density" sheets, and switch
# the active sheet/selection to the "Medium traffic density" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Low traffic densit" -----------------------------------
$ws1 = $wb.Worksheets.Item("Low traffic densit")

$ws1.Range("C18").Value = 24
$ws1.Range("D18").Value = 4561.95
$ws1.Range("E18").Value = 515.3

$ws1.Range("C19").Value = 18
$ws1.Range("D19").Value = 4878.58
$ws1.Range("E19").Value = 443.61

$ws1.Range("C20").Value = 20
$ws1.Range("D20").Value = 4652.38
$ws1.Range("E20").Value = 494.49

$ws1.Range("C21").Value = 24
$ws1.Range("D21").Value = 4578.98
$ws1.Range("E21").Value = 498.17

$ws1.Range("C22").Value = 21
$ws1.Range("D22").Value = 4682.81
$ws1.Range("E22").Value = 489.22

$ws1.Range("C23").Value = 20
$ws1.Range("D23").Value = 4540.08
$ws1.Range("E23").Value = 460.51

$ws1.Range("C24").Value = 21
$ws1.Range("D24").Value = 4374.59
$ws1.Range("E24").Value = 466.45

$ws1.Range("C25").Value = 20
$ws1.Range("D25").Value = 4950.2
$ws1.Range("E25").Value = 511.96

$ws1.Range("C26").Value = 25
$ws1.Range("D26").Value = 4885.72
$ws1.Range("E26").Value = 538.26

# --- Sheet 2: "Medium traffic density" -------------------------------
$ws2 = $wb.Worksheets.Item("Medium traffic density")

$ws2.Range("C17").Value = 130
$ws2.Range("D17").Value = 3937.23
$ws2.Range("E17").Value = 491.33

$ws2.Range("C18").Value = 121
$ws2.Range("D18").Value = 3233.98
$ws2.Range("E18").Value = 431.47

$ws2.Range("C19").Value = 135
$ws2.Range("D19").Value = 3231.53
$ws2.Range("E19").Value = 416.29

$ws2.Range("C20").Value = 127
$ws2.Range("D20").Value = 3397.79
$ws2.Range("E20").Value = 445.49

$ws2.Range("C21").Value = 139
$ws2.Range("D21").Value = 3193.53
$ws2.Range("E21").Value = 414.29

$ws2.Range("C22").Value = 124
$ws2.Range("D22").Value = 3436.9
$ws2.Range("E22").Value = 450.48

$ws2.Range("C23").Value = 131
$ws2.Range("D23").Value = 3516.81
$ws2.Range("E23").Value = 485.36

$ws2.Range("C24").Value = 132
$ws2.Range("D24").Value = 3598.33
$ws2.Range("E24").Value = 491.52

$ws2.Range("C25").Value = 138
$ws2.Range("D25").Value = 3329.31
$ws2.Range("E25").Value = 463.21

$ws2.Range("C26").Value = 134
$ws2.Range("D26").Value = 3305.63
$ws2.Range("E26").Value = 443.36

# --- View / selection state -------------------------------------------
# "Low traffic densit" becomes non-selected, scrolled to A4, selection D19
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("D19").Select()

# "Medium traffic density" becomes the active/selected tab, scrolled to
# A4, selection C18
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("C18").Select()
